$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 38.33049633333334
$ws.Range("H2").Value = 114.991489
$ws.Range("I2").Value = 0.5317874798120843
$ws.Range("J2").Value = 0.5317874798120843
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003986333333333333
$ws.Range("N2").Value = 0.011959
$ws.Range("O2").Value = 0.0001042515924716245
$ws.Range("P2").Value = 0.0001042515924716245
$ws.Range("Q2").Value = 0.1527981352167778
$ws.Range("R2").Value = 1.375183216951
$ws.Range("S2").Value = 0.00005543969162688167
$ws.Range("T2").Value = 0.00005543969162688167
$ws.Range("G3").Value = 38.33049633333334
$ws.Range("H3").Value = 114.991489
$ws.Range("I3").Value = 0.5317874798120843
$ws.Range("J3").Value = 0.5317874798120843
$ws.Range("O3").Value = 0.005104344171442446
$ws.Range("P3").Value = 0.005104344171442446
$ws.Range("Q3").Value = 7.481269613347333
$ws.Range("R3").Value = 67.331426520126
$ws.Range("S3").Value = 0.00271442632302488
$ws.Range("T3").Value = 0.00271442632302488
$ws.Range("G4").Value = 38.33049633333334
$ws.Range("H4").Value = 114.991489
$ws.Range("I4").Value = 0.5317874798120843
$ws.Range("J4").Value = 0.5317874798120843
$ws.Range("M4").Value = 20.492743
$ws.Range("N4").Value = 61.478229
$ws.Range("O4").Value = 0.5359313718191496
$ws.Range("P4").Value = 0.5359313718191496
$ws.Range("Q4").Value = 785.4970104214424
$ws.Range("R4").Value = 7069.473093792981
$ws.Range("S4").Value = 0.2850015935719387
$ws.Range("T4").Value = 0.2850015935719387
$ws.Range("G5").Value = 38.33049633333334
$ws.Range("H5").Value = 114.991489
$ws.Range("I5").Value = 0.5317874798120843
$ws.Range("J5").Value = 0.5317874798120843
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1451653333333333
$ws.Range("N5").Value = 0.435496
$ws.Range("O5").Value = 0.003796400327370399
$ws.Range("P5").Value = 0.003796400327370399
$ws.Range("Q5").Value = 5.564259277060446
$ws.Range("R5").Value = 50.078333493544
$ws.Range("S5").Value = 0.002018878162450076
$ws.Range("T5").Value = 0.002018878162450076
$ws.Range("G6").Value = 38.33049633333334
$ws.Range("H6").Value = 114.991489
$ws.Range("I6").Value = 0.5317874798120843
$ws.Range("J6").Value = 0.5317874798120843
$ws.Range("M6").Value = 17.40055266666667
$ws.Range("N6").Value = 52.20165799999999
$ws.Range("O6").Value = 0.4550636320895659
$ws.Range("P6").Value = 0.4550636320895659
$ws.Range("Q6").Value = 666.9718201876402
$ws.Range("R6").Value = 6002.746381688761
$ws.Range("S6").Value = 0.2419971420630438
$ws.Range("T6").Value = 0.2419971420630438
$ws.Range("I7").Value = 0.2073226210890634
$ws.Range("J7").Value = 0.2073226210890634
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.003986333333333333
$ws.Range("N7").Value = 0.011959
$ws.Range("O7").Value = 0.0001042515924716245
$ws.Range("P7").Value = 0.0001042515924716245
$ws.Range("Q7").Value = 0.05956986783866667
$ws.Range("R7").Value = 0.536128810548
$ws.Range("S7").Value = 0.00002161371340392606
$ws.Range("T7").Value = 0.00002161371340392606
$ws.Range("I8").Value = 0.2073226210890634
$ws.Range("J8").Value = 0.2073226210890634
$ws.Range("O8").Value = 0.005104344171442446
$ws.Range("P8").Value = 0.005104344171442446
$ws.Range("S8").Value = 0.001058246012564131
$ws.Range("T8").Value = 0.001058246012564131
$ws.Range("I9").Value = 0.2073226210890634
$ws.Range("J9").Value = 0.2073226210890634
$ws.Range("M9").Value = 20.492743
$ws.Range("N9").Value = 61.478229
$ws.Range("O9").Value = 0.5359313718191496
$ws.Range("P9").Value = 0.5359313718191496
$ws.Range("Q9").Value = 306.233796846332
$ws.Range("R9").Value = 2756.104171616988
$ws.Range("S9").Value = 0.1111106967294035
$ws.Range("T9").Value = 0.1111106967294035
$ws.Range("I10").Value = 0.2073226210890634
$ws.Range("J10").Value = 0.2073226210890634
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1451653333333333
$ws.Range("N10").Value = 0.435496
$ws.Range("O10").Value = 0.003796400327370399
$ws.Range("P10").Value = 0.003796400327370399
$ws.Range("Q10").Value = 2.169281642634667
$ws.Range("R10").Value = 19.523534783712
$ws.Range("S10").Value = 0.0007870796665738094
$ws.Range("T10").Value = 0.0007870796665738094
$ws.Range("I11").Value = 0.2073226210890634
$ws.Range("J11").Value = 0.2073226210890634
$ws.Range("M11").Value = 17.40055266666667
$ws.Range("N11").Value = 52.20165799999999
$ws.Range("O11").Value = 0.4550636320895659
$ws.Range("P11").Value = 0.4550636320895659
$ws.Range("Q11").Value = 260.0255763875974
$ws.Range("R11").Value = 2340.230187488376
$ws.Range("S11").Value = 0.094344984967118
$ws.Range("T11").Value = 0.094344984967118
$ws.Range("G12").Value = 8.167063666666666
$ws.Range("H12").Value = 24.501191
$ws.Range("I12").Value = 0.1133077476219524
$ws.Range("J12").Value = 0.1133077476219524
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.003986333333333333
$ws.Range("N12").Value = 0.011959
$ws.Range("O12").Value = 0.0001042515924716245
$ws.Range("P12").Value = 0.0001042515924716245
$ws.Range("Q12").Value = 0.03255663812988888
$ws.Range("R12").Value = 0.293009743169
$ws.Range("S12").Value = 0.00001181251312896147
$ws.Range("T12").Value = 0.00001181251312896147
$ws.Range("G13").Value = 8.167063666666666
$ws.Range("H13").Value = 24.501191
$ws.Range("I13").Value = 0.1133077476219524
$ws.Range("J13").Value = 0.1133077476219524
$ws.Range("O13").Value = 0.005104344171442446
$ws.Range("P13").Value = 0.005104344171442446
$ws.Range("Q13").Value = 1.594031152332666
$ws.Range("R13").Value = 14.346280370994
$ws.Range("S13").Value = 0.0005783617411533845
$ws.Range("T13").Value = 0.0005783617411533847
$ws.Range("G14").Value = 8.167063666666666
$ws.Range("H14").Value = 24.501191
$ws.Range("I14").Value = 0.1133077476219524
$ws.Range("J14").Value = 0.1133077476219524
$ws.Range("M14").Value = 20.492743
$ws.Range("N14").Value = 61.478229
$ws.Range("O14").Value = 0.5359313718191496
$ws.Range("P14").Value = 0.5359313718191496
$ws.Range("Q14").Value = 167.3655367856377
$ws.Range("R14").Value = 1506.289831070739
$ws.Range("S14").Value = 0.06072517662077095
$ws.Range("T14").Value = 0.06072517662077096
$ws.Range("G15").Value = 8.167063666666666
$ws.Range("H15").Value = 24.501191
$ws.Range("I15").Value = 0.1133077476219524
$ws.Range("J15").Value = 0.1133077476219524
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.1451653333333333
$ws.Range("N15").Value = 0.435496
$ws.Range("O15").Value = 0.003796400327370399
$ws.Range("P15").Value = 0.003796400327370399
$ws.Range("Q15").Value = 1.185574519526222
$ws.Range("R15").Value = 10.670170675736
$ws.Range("S15").Value = 0.0004301615701655828
$ws.Range("T15").Value = 0.0004301615701655829
$ws.Range("G16").Value = 8.167063666666666
$ws.Range("H16").Value = 24.501191
$ws.Range("I16").Value = 0.1133077476219524
$ws.Range("J16").Value = 0.1133077476219524
$ws.Range("M16").Value = 17.40055266666667
$ws.Range("N16").Value = 52.20165799999999
$ws.Range("O16").Value = 0.4550636320895659
$ws.Range("P16").Value = 0.4550636320895659
$ws.Range("Q16").Value = 142.1114214638531
$ws.Range("R16").Value = 1279.002793174678
$ws.Range("S16").Value = 0.05156223517673355
$ws.Range("T16").Value = 0.05156223517673356
$ws.Range("G17").Value = 5.834252333333334
$ws.Range("H17").Value = 17.502757
$ws.Range("I17").Value = 0.08094292121735479
$ws.Range("J17").Value = 0.08094292121735479
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.003986333333333333
$ws.Range("N17").Value = 0.011959
$ws.Range("O17").Value = 0.0001042515924716245
$ws.Range("P17").Value = 0.0001042515924716245
$ws.Range("Q17").Value = 0.02325727455144445
$ws.Range("R17").Value = 0.209315470963
$ws.Range("S17").Value = 0.000008438428436214482
$ws.Range("T17").Value = 0.000008438428436214482
$ws.Range("G18").Value = 5.834252333333334
$ws.Range("H18").Value = 17.502757
$ws.Range("I18").Value = 0.08094292121735479
$ws.Range("J18").Value = 0.08094292121735479
$ws.Range("O18").Value = 0.005104344171442446
$ws.Range("P18").Value = 0.005104344171442446
$ws.Range("Q18").Value = 1.138717701915333
$ws.Range("R18").Value = 10.248459317238
$ws.Range("S18").Value = 0.00041316052813533
$ws.Range("T18").Value = 0.00041316052813533
$ws.Range("G19").Value = 5.834252333333334
$ws.Range("H19").Value = 17.502757
$ws.Range("I19").Value = 0.08094292121735479
$ws.Range("J19").Value = 0.08094292121735479
$ws.Range("M19").Value = 20.492743
$ws.Range("N19").Value = 61.478229
$ws.Range("O19").Value = 0.5359313718191496
$ws.Range("P19").Value = 0.5359313718191496
$ws.Range("Q19").Value = 119.5598336641504
$ws.Range("R19").Value = 1076.038502977353
$ws.Range("S19").Value = 0.04337985080706631
$ws.Range("T19").Value = 0.04337985080706631
$ws.Range("G20").Value = 5.834252333333334
$ws.Range("H20").Value = 17.502757
$ws.Range("I20").Value = 0.08094292121735479
$ws.Range("J20").Value = 0.08094292121735479
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.1451653333333333
$ws.Range("N20").Value = 0.435496
$ws.Range("O20").Value = 0.003796400327370399
$ws.Range("P20").Value = 0.003796400327370399
$ws.Range("Q20").Value = 0.8469311847191113
$ws.Range("R20").Value = 7.622380662472001
$ws.Range("S20").Value = 0.0003072917326078821
$ws.Range("T20").Value = 0.0003072917326078821
$ws.Range("G21").Value = 5.834252333333334
$ws.Range("H21").Value = 17.502757
$ws.Range("I21").Value = 0.08094292121735479
$ws.Range("J21").Value = 0.08094292121735479
$ws.Range("M21").Value = 17.40055266666667
$ws.Range("N21").Value = 52.20165799999999
$ws.Range("O21").Value = 0.4550636320895659
$ws.Range("P21").Value = 0.4550636320895659
$ws.Range("Q21").Value = 101.5192149967896
$ws.Range("R21").Value = 913.6729349711061
$ws.Range("S21").Value = 0.03683417972110906
$ws.Range("T21").Value = 0.03683417972110906
$ws.Range("G22").Value = 4.803262333333334
$ws.Range("H22").Value = 14.409787
$ws.Range("I22").Value = 0.06663923025954499
$ws.Range("J22").Value = 0.066639230259545
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0.3333333333333333
$ws.Range("M22").Value = 0.003986333333333333
$ws.Range("N22").Value = 0.011959
$ws.Range("O22").Value = 0.0001042515924716245
$ws.Range("P22").Value = 0.0001042515924716245
$ws.Range("Q22").Value = 0.01914740474811111
$ws.Range("R22").Value = 0.172326642733
$ws.Range("S22").Value = 0.000006947245875640834
$ws.Range("T22").Value = 0.000006947245875640836
$ws.Range("G23").Value = 4.803262333333334
$ws.Range("H23").Value = 14.409787
$ws.Range("I23").Value = 0.06663923025954499
$ws.Range("J23").Value = 0.066639230259545
$ws.Range("O23").Value = 0.005104344171442446
$ws.Range("P23").Value = 0.005104344171442446
$ws.Range("Q23").Value = 0.9374911356953334
$ws.Range("R23").Value = 8.437420221258002
$ws.Range("S23").Value = 0.0003401495665647195
$ws.Range("T23").Value = 0.0003401495665647196
$ws.Range("G24").Value = 4.803262333333334
$ws.Range("H24").Value = 14.409787
$ws.Range("I24").Value = 0.06663923025954499
$ws.Range("J24").Value = 0.066639230259545
$ws.Range("M24").Value = 20.492743
$ws.Range("N24").Value = 61.478229
$ws.Range("O24").Value = 0.5359313718191496
$ws.Range("P24").Value = 0.5359313718191496
$ws.Range("Q24").Value = 98.43202055858035
$ws.Range("R24").Value = 885.8881850272231
$ws.Range("S24").Value = 0.03571405408997013
$ws.Range("T24").Value = 0.03571405408997014
$ws.Range("G25").Value = 4.803262333333334
$ws.Range("H25").Value = 14.409787
$ws.Range("I25").Value = 0.06663923025954499
$ws.Range("J25").Value = 0.066639230259545
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.1451653333333333
$ws.Range("N25").Value = 0.435496
$ws.Range("O25").Value = 0.003796400327370399
$ws.Range("P25").Value = 0.003796400327370399
$ws.Range("Q25").Value = 0.6972671777057778
$ws.Range("R25").Value = 6.275404599352001
$ws.Range("S25").Value = 0.000252989195573048
$ws.Range("T25").Value = 0.0002529891955730481
$ws.Range("G26").Value = 4.803262333333334
$ws.Range("H26").Value = 14.409787
$ws.Range("I26").Value = 0.06663923025954499
$ws.Range("J26").Value = 0.066639230259545
$ws.Range("M26").Value = 17.40055266666667
$ws.Range("N26").Value = 52.20165799999999
$ws.Range("O26").Value = 0.4550636320895659
$ws.Range("P26").Value = 0.4550636320895659
$ws.Range("Q26").Value = 83.57941920298289
$ws.Range("R26").Value = 752.214772826846
$ws.Range("S26").Value = 0.03032509016156145
$ws.Range("T26").Value = 0.03032509016156145
